$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells G1:J1 ---
$ws.Range("G1").Value = "stageId"
$ws.Range("H1").Value = "stageName"
$ws.Range("I1").Value = "slot"
$ws.Range("J1").Value = "datetime"

# Copy the header formatting (bold/centered/bordered) from the existing
# header row onto the new header cells so they reuse the same style.
$ws.Range("A1").Copy()
$ws.Range("G1:J1").PasteSpecial(-4122)

# --- Row 2 fix-ups ---
# E2 was stored as a text value of the phone number; it becomes numeric.
$ws.Range("E2").Value = 1234567890

$ws.Range("G2").Value = 1
$ws.Range("H2").Value = "A"
$ws.Range("I2").Value = "Slot 1"
$ws.Range("J2").Value = "2025-01-29 20:23:18.351679"

# --- New row 3 ---
$ws.Range("A3").Value = "w"
$ws.Range("B3").Value = "w"
$ws.Range("C3").Value = "w"
$ws.Range("D3").Value = "w"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = "A"
$ws.Range("I3").Value = "Slot 1"
$ws.Range("J3").Value = "2025-01-29 20:25:44.348060"

# --- New row 4 ---
$ws.Range("A4").Value = "e"
$ws.Range("B4").Value = "e"
$ws.Range("C4").Value = "e"
$ws.Range("D4").Value = "e"

# E4 must stay a text value "4" (not the number 4). Stage the text via a
# scratch cell (quote-prefixed so it is stored as text) and paste only the
# value across, so E4's formatting/style is left untouched.
$ws.Range("Z1").Value = "'4"
$ws.Range("Z1").Copy()
$ws.Range("E4").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = "A"
$ws.Range("I4").Value = "Slot 2"
$ws.Range("J4").Value = "2025-01-29 20:26:03.309266"
